{"js": "// Update the \"interest\" requirement row:\n//  - \"...allow accounts to collect interest\" -> \"...allow savings accounts to collect interest\"\n//  - \"...flat 0.5% APY\" -> \"...flat 0.05% APY\"\n\nconst firstPart = context.document.body.search(\n  \"This system shall allow accounts to collect interest\",\n  { matchCase: true }\n);\nfirstPart.load(\"items\");\nawait context.sync();\n\nif (firstPart.items.length > 0) {\n  firstPart.items[0].insertText(\n    \"This system shall allow savings accounts to collect interest\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\nconst secondPart = context.document.body.search(\n  \"calculated as a flat 0.5% APY\",\n  { matchCase: true }\n);\nsecondPart.load(\"items\");\nawait context.sync();\n\nif (secondPart.items.length > 0) {\n  secondPart.items[0].insertText(\n    \"calculated as a flat 0.05% APY\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n", "ps1": "# Update the \"interest\" requirement row:\n#  - \"...allow accounts to collect interest\" -> \"...allow savings accounts to collect interest\"\n#  - \"...flat 0.5% APY\" -> \"...flat 0.05% APY\"\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"This system shall allow accounts to collect interest\"\n$find1.Replacement.Text = \"This system shall allow savings accounts to collect interest\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find1.Replacement.Text, $wdReplaceAll)\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"calculated as a flat 0.5% APY\"\n$find2.Replacement.Text = \"calculated as a flat 0.05% APY\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find2.Replacement.Text, $wdReplaceAll)\n"}
